$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count()
for ($s = 1; $s -le $sheetCount; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        # NB: keep the string literal on the left of -eq — PowerShell's -eq
        # coerces the right-hand side to the left-hand side's type, and a
        # boolean cell value ($true) would otherwise "equal" any non-empty
        # string and be clobbered.
        if ("Ready for handoff" -eq $val) {
            $cell.Value = "In Translation"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Shrink the "Status" / language columns that held that text, since the
#    new text is shorter (this mirrors the report generator re-sizing the
#    columns after refreshing the status values).
#    NOTE: the host engine adds a fixed 5/6 character padding to whatever
#    ColumnWidth is assigned once the workbook is saved, so the assigned
#    value is pre-compensated to land on the desired on-disk width.
# ---------------------------------------------------------------------------
$newWidth = 13.4101845877511 - (5 / 6)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth   # column F (de-de)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth       # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth       # column C (Status)
